$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = [char]39

function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).Value2 = "$apos$val"
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws 'D2' '41.408.23'
Set-TextValue $ws 'E2' '  -1.21%  '
Set-TextValue $ws 'D3' '2.158.47'
Set-TextValue $ws 'E3' '  -2.97%  '
Set-TextValue $ws 'E4' '  +0.39%  '
Set-TextValue $ws 'D5' '237.18'
Set-TextValue $ws 'E5' '  -1.97%  '
Set-TextValue $ws 'E6' '  -3.46%  '
Set-TextValue $ws 'D7' '71.13'
Set-TextValue $ws 'E7' '  -2.55%  '
Set-TextValue $ws 'E8' '  +0.08%  '
Set-TextValue $ws 'D9' '0.572'
Set-TextValue $ws 'E9' '  -4.47%  '
Set-TextValue $ws 'D10' '39.46'
Set-TextValue $ws 'E10' '  -7.10%  '
Set-TextValue $ws 'D11' '0.0900'
Set-TextValue $ws 'E11' '  -5.62%  '
Set-TextValue $ws 'D12' '53.90'
Set-TextValue $ws 'E12' '  -4.84%  '
Set-TextValue $ws 'D13' '0.0998'
Set-TextValue $ws 'E13' '  -3.80%  '
Set-TextValue $ws 'D14' '6.65'
Set-TextValue $ws 'E14' '  -4.69%  '
Set-TextValue $ws 'D15' '2.482.55'
Set-TextValue $ws 'E15' '  -2.87%  '
Set-TextValue $ws 'D16' '14.09'
Set-TextValue $ws 'E16' '  -1.48%  '
Set-TextValue $ws 'D17' '2.147.88'
Set-TextValue $ws 'E17' '  -2.77%  '
Set-TextValue $ws 'D18' '0.777'
Set-TextValue $ws 'E18' '  -7.24%  '
Set-TextValue $ws 'D19' '41.340.19'
Set-TextValue $ws 'E19' '  -0.99%  '
Set-TextValue $ws 'E20' '  -4.50%  '
Set-TextValue $ws 'D21' '69.48'
Set-TextValue $ws 'E21' '  -4.41%  '
Set-TextValue $ws 'E22' '  -6.96%  '
Set-TextValue $ws 'D23' '9.91'
Set-TextValue $ws 'E23' '  -12.19%  '
Set-TextValue $ws 'D24' '226.86'
Set-TextValue $ws 'E24' '  -1.29%  '
Set-TextValue $ws 'E25' '  -4.43%  '
Set-TextValue $ws 'E26' '  -0.15%  '
Set-TextValue $ws 'D27' '10.61'
Set-TextValue $ws 'E27' '  -7.03%  '
Set-TextValue $ws 'D28' '3.28'
Set-TextValue $ws 'E28' '  -9.45%  '
Set-TextValue $ws 'D29' '2.18'
Set-TextValue $ws 'E29' '  -4.62%  '
Set-TextValue $ws 'E30' '  -0.93%  '
Set-TextValue $ws 'D31' '171.40'
Set-TextValue $ws 'E31' '  +2.65%  '
Set-TextValue $ws 'D32' '19.69'
Set-TextValue $ws 'E32' '  -3.95%  '
Set-TextValue $ws 'D33' '32.82'
Set-TextValue $ws 'E33' '  +8.15%  '
Set-TextValue $ws 'D34' '0.0765'
Set-TextValue $ws 'E34' '  -4.22%  '
Set-TextValue $ws 'D35' '5.09'
Set-TextValue $ws 'E35' '  -8.33%  '
Set-TextValue $ws 'E36' '  -4.08%  '
Set-TextValue $ws 'E37' '  -4.17%  '
Set-TextValue $ws 'D38' '4.20'
Set-TextValue $ws 'E38' '  -2.11%  '
Set-TextValue $ws 'D39' '0.0297'
Set-TextValue $ws 'E39' '  -2.09%  '
Set-TextValue $ws 'D40' '11.87'
Set-TextValue $ws 'E40' '  -11.34%  '
Set-TextValue $ws 'E41' '  -3.15%  '
Set-TextValue $ws 'D42' '5.30'
Set-TextValue $ws 'E42' '  -6.23%  '
Set-TextValue $ws 'D43' '58.39'
Set-TextValue $ws 'E43' '  -10.50%  '
Set-TextValue $ws 'E44' '  -5.10%  '
Set-TextValue $ws 'D45' '8.36'
Set-TextValue $ws 'E45' '  -4.31%  '
Set-TextValue $ws 'D46' '0.0953'
Set-TextValue $ws 'E46' '  -5.25%  '
Set-TextValue $ws 'D47' '95.84'
Set-TextValue $ws 'E47' '  -7.71%  '
Set-TextValue $ws 'E48' '  -3.89%  '
Set-TextValue $ws 'E49' '  -5.29%  '
Set-TextValue $ws 'B50' 'HuobiToken'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D50' '2.62'
Set-TextValue $ws 'E50' '  -2.63%  '
Set-TextValue $ws 'B51' 'NEARProtocol'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D51' '2.15'
Set-TextValue $ws 'E51' '  -8.08%  '
